# "break out stock.yaml completed"
# Fill in the header row + first data row on the "3 V 0.3" sheet, matching
# the other already-populated breakout sheets (sheet1 "10per change",
# sheet3 "DND 3 V 0.3").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("3 V 0.3")

# --- Header row (row 1) -----------------------------------------------
# A1 already holds "Date Time" with the bold/centered header style;
# extend the same header labels across B1:H1 and copy A1's formatting
# onto them so they match (bold, centered, thin border).
$ws.Range("B1").Value = "sr"
$ws.Range("C1").Value = "nsecode"
$ws.Range("D1").Value = "name"
$ws.Range("E1").Value = "bsecode"
$ws.Range("F1").Value = "per_chg"
$ws.Range("G1").Value = "close"
$ws.Range("H1").Value = "volume"

$ws.Range("A1").Copy()
$ws.Range("B1:H1").PasteSpecial(-4122)   # xlPasteFormats

# --- Data row (row 2) ---------------------------------------------------
$ws.Range("A2").Value = "12/06/2024 05:44:55"
$ws.Range("B2").Value = 1
$ws.Range("C2").Value = "PAISALO"
$ws.Range("D2").Value = "Paisalo Digital Ltd"

# bsecode on this row is stored as text ("532900"), not a number - force
# text entry without leaving a lingering custom number format behind.
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "532900"
$ws.Range("E2").Style = "Normal"

$ws.Range("F2").Value = 8.24
$ws.Range("G2").Value = 68.8
$ws.Range("H2").Value = 2503862
